# Generate Report for handoff
#
# - Status moves from "Handoff transform failed" to "Ready for handoff"
#   (Overview summary + both language detail sheets).
# - zh-cn / de-de detail sheets each get their "Latest Handoff File" (C2)
#   populated with a hyperlink to the generated .xlf handoff file, the
#   "Latest Handoff Datetime" (D2) stamped, and the "Handoff Reason" (H2)
#   switched from "Ignored" to "Include" now that a handoff went out.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# --- Overview sheet: both language status cells reflect the new status ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus

# --- zh-cn detail sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B2").Value = $newStatus

$zhcnFile = "aaad06e6-4e31-498a-8189-d929387b967f.54d9a52e53430d0f502f323ed4fb59292021e391.zh-cn.xlf"
$zhcnUrl = "https://github.com/OpenLocalizationTest/oltest/blob/4da00e306447362d65ca424a631e17b88bd7ccd4/e2e/" + $zhcnFile
$zhcn.Hyperlinks.Add($zhcn.Range("C2"), $zhcnUrl, "", "", $zhcnFile)

$zhcn.Range("D2").Value = "2016-01-11 17:14:23"
$zhcn.Range("H2").Value = "Include"

# --- de-de detail sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B2").Value = $newStatus

$dedeFile = "aaad06e6-4e31-498a-8189-d929387b967f.54d9a52e53430d0f502f323ed4fb59292021e391.de-de.xlf"
$dedeUrl = "https://github.com/OpenLocalizationTest/oltest/blob/4da00e306447362d65ca424a631e17b88bd7ccd4/e2e/" + $dedeFile
$dede.Hyperlinks.Add($dede.Range("C2"), $dedeUrl, "", "", $dedeFile)

$dede.Range("D2").Value = "2016-01-11 17:14:43"
$dede.Range("H2").Value = "Include"

Write-Output "Report generated for handoff"
